# Update database and change read_price algorithm
# Shift the yearly income-statement columns (D:H) one period to the left,
# dropping the oldest period (1396/10) and appending the new period
# (1401/10) together with its corresponding values / publish dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(4, 5, 6, 7, 8)   # D, E, F, G, H

# --- Row 8: financial period headers -------------------------------------
$period8 = @(
    "12 ماهه منتهی به 1397/10",
    "12 ماهه منتهی به 1398/10",
    "12 ماهه منتهی به 1399/10",
    "12 ماهه منتهی به 1400/10",
    "12 ماهه منتهی به 1401/10"
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(8, $cols[$i]).Value2 = $period8[$i]
}

# --- Row 9: publish dates --------------------------------------------------
$dates9 = @(
    "1399-04-02 (10)",
    "1400-02-26 (8)",
    "1401-02-19 (9)",
    "1402-02-13 (9)",
    "1402-02-29 (3)"
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(9, $cols[$i]).Value2 = $dates9[$i]
}

# --- Numeric data rows: shift left and append new last column -------------
$rowValues = @{
    11 = @(16042, 24117, 19903, 24119, 32466)
    12 = @(-10812, -12273, -8843, -11682, -15216)
    13 = @(5230, 11844, 11060, 12437, 17250)
    14 = @(-1662, -4339, -2724, -2960, -4321)
    16 = @(121, -998, -17, -387, -43)
    17 = @(3689, 6507, 8318, 9089, 12887)
    18 = @(-1895, -1272, -565, -256, -48)
    19 = @(115, 565, 1237, 3214, 3380)
    20 = @(1909, 5799, 8990, 12048, 16218)
    21 = @(-172, -943, -929, -1569, -1342)
    22 = @(1737, 4857, 8062, 10479, 14876)
    24 = @(1737, 4857, 8062, 10479, 14876)
    26 = @(5650, 5820, 3479, 2809, 2297)
}

foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value2 = $vals[$i]
    }
}
